# Fix formatting when scraping floating point numbers:
# - Replace comma-decimal "Importe" text values (e.g. "3.398,50") with dot-decimal ("3398.50")
#   in column H (rows 2-240), without altering the existing cell style.
# - Replace stray commas used as name separators with periods in a handful of
#   "Razon social"/"Nombre Fantasia" cells, and normalize "S.H." -> "SH".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H: re-format amounts, keep default (unstyled) text cells ---
$hRange = $ws.Range("H2:H240")
$hRange.NumberFormat = "@"
$hValues = New-Object 'object[,]' 239,1
$hValues[0,0] = '3398.50'
$hValues[1,0] = '35400.00'
$hValues[2,0] = '2350.00'
$hValues[3,0] = '853.98'
$hValues[4,0] = '9750.00'
$hValues[5,0] = '35980.88'
$hValues[6,0] = '199500.00'
$hValues[7,0] = '119000.00'
$hValues[8,0] = '1509.70'
$hValues[9,0] = '1200.00'
$hValues[10,0] = '1682.20'
$hValues[11,0] = '90.00'
$hValues[12,0] = '1700.00'
$hValues[13,0] = '108351.10'
$hValues[14,0] = '107242.98'
$hValues[15,0] = '6530.00'
$hValues[16,0] = '1124.20'
$hValues[17,0] = '60722.63'
$hValues[18,0] = '3521.13'
$hValues[19,0] = '2789.89'
$hValues[20,0] = '9139.23'
$hValues[21,0] = '14232.10'
$hValues[22,0] = '2000.00'
$hValues[23,0] = '99.28'
$hValues[24,0] = '1810.00'
$hValues[25,0] = '13.80'
$hValues[26,0] = '22100.00'
$hValues[27,0] = '29775.56'
$hValues[28,0] = '26.48'
$hValues[29,0] = '174.50'
$hValues[30,0] = '4201.94'
$hValues[31,0] = '100.05'
$hValues[32,0] = '3524.00'
$hValues[33,0] = '10400.00'
$hValues[34,0] = '30.00'
$hValues[35,0] = '12.80'
$hValues[36,0] = '255.25'
$hValues[37,0] = '12767.90'
$hValues[38,0] = '15600.00'
$hValues[39,0] = '521.61'
$hValues[40,0] = '734.13'
$hValues[41,0] = '590.00'
$hValues[42,0] = '175.00'
$hValues[43,0] = '150983.41'
$hValues[44,0] = '2321.85'
$hValues[45,0] = '279.00'
$hValues[46,0] = '838.60'
$hValues[47,0] = '2780.83'
$hValues[48,0] = '1120.00'
$hValues[49,0] = '2220.44'
$hValues[50,0] = '21.30'
$hValues[51,0] = '5660.81'
$hValues[52,0] = '26054.28'
$hValues[53,0] = '39.37'
$hValues[54,0] = '1360.00'
$hValues[55,0] = '1458.12'
$hValues[56,0] = '948.00'
$hValues[57,0] = '4337.00'
$hValues[58,0] = '469.00'
$hValues[59,0] = '21736.00'
$hValues[60,0] = '72.00'
$hValues[61,0] = '2914.80'
$hValues[62,0] = '731.36'
$hValues[63,0] = '1181.75'
$hValues[64,0] = '101065.22'
$hValues[65,0] = '2674.00'
$hValues[66,0] = '1991.69'
$hValues[67,0] = '975.00'
$hValues[68,0] = '399.00'
$hValues[69,0] = '596.03'
$hValues[70,0] = '2714.00'
$hValues[71,0] = '2867.01'
$hValues[72,0] = '9697.00'
$hValues[73,0] = '16667.00'
$hValues[74,0] = '1796.80'
$hValues[75,0] = '465.00'
$hValues[76,0] = '960.00'
$hValues[77,0] = '75.00'
$hValues[78,0] = '60.00'
$hValues[79,0] = '11900.00'
$hValues[80,0] = '440.00'
$hValues[81,0] = '2176.80'
$hValues[82,0] = '2239.00'
$hValues[83,0] = '249.00'
$hValues[84,0] = '1205.00'
$hValues[85,0] = '320.00'
$hValues[86,0] = '195.33'
$hValues[87,0] = '6065.00'
$hValues[88,0] = '386.65'
$hValues[89,0] = '545.08'
$hValues[90,0] = '1309.00'
$hValues[91,0] = '168.00'
$hValues[92,0] = '1600.00'
$hValues[93,0] = '165.00'
$hValues[94,0] = '310.00'
$hValues[95,0] = '608.40'
$hValues[96,0] = '5000.00'
$hValues[97,0] = '116893.20'
$hValues[98,0] = '3.12'
$hValues[99,0] = '116.16'
$hValues[100,0] = '52.15'
$hValues[101,0] = '124.14'
$hValues[102,0] = '690.00'
$hValues[103,0] = '6045.09'
$hValues[104,0] = '489.59'
$hValues[105,0] = '29250.00'
$hValues[106,0] = '900.00'
$hValues[107,0] = '4895.69'
$hValues[108,0] = '3861.50'
$hValues[109,0] = '200.00'
$hValues[110,0] = '875.00'
$hValues[111,0] = '517.50'
$hValues[112,0] = '120.00'
$hValues[113,0] = '220.80'
$hValues[114,0] = '180.00'
$hValues[115,0] = '436.50'
$hValues[116,0] = '1148.05'
$hValues[117,0] = '779.00'
$hValues[118,0] = '20.67'
$hValues[119,0] = '4677.20'
$hValues[120,0] = '4895.08'
$hValues[121,0] = '568.00'
$hValues[122,0] = '550.00'
$hValues[123,0] = '3014.00'
$hValues[124,0] = '1000.00'
$hValues[125,0] = '2375.00'
$hValues[126,0] = '760.00'
$hValues[127,0] = '2926.50'
$hValues[128,0] = '6080.00'
$hValues[129,0] = '440.00'
$hValues[130,0] = '1000.00'
$hValues[131,0] = '5800.00'
$hValues[132,0] = '2800.00'
$hValues[133,0] = '6000.00'
$hValues[134,0] = '4300.00'
$hValues[135,0] = '5350.00'
$hValues[136,0] = '3609.90'
$hValues[137,0] = '1332.25'
$hValues[138,0] = '300.00'
$hValues[139,0] = '635.00'
$hValues[140,0] = '33.00'
$hValues[141,0] = '4951.00'
$hValues[142,0] = '292.68'
$hValues[143,0] = '134800.00'
$hValues[144,0] = '2535.00'
$hValues[145,0] = '28.00'
$hValues[146,0] = '132810.00'
$hValues[147,0] = '5430.00'
$hValues[148,0] = '1600.00'
$hValues[149,0] = '1000.00'
$hValues[150,0] = '6806.25'
$hValues[151,0] = '384.00'
$hValues[152,0] = '600.00'
$hValues[153,0] = '1000.00'
$hValues[154,0] = '10588.50'
$hValues[155,0] = '1500.00'
$hValues[156,0] = '950.00'
$hValues[157,0] = '750.00'
$hValues[158,0] = '2550.79'
$hValues[159,0] = '1500.00'
$hValues[160,0] = '200.00'
$hValues[161,0] = '290.00'
$hValues[162,0] = '9440.00'
$hValues[163,0] = '2000.00'
$hValues[164,0] = '50.00'
$hValues[165,0] = '30.00'
$hValues[166,0] = '560.00'
$hValues[167,0] = '2363.00'
$hValues[168,0] = '550.00'
$hValues[169,0] = '50.00'
$hValues[170,0] = '2030.00'
$hValues[171,0] = '174.59'
$hValues[172,0] = '2577.96'
$hValues[173,0] = '480.00'
$hValues[174,0] = '137.20'
$hValues[175,0] = '1468.00'
$hValues[176,0] = '2770.00'
$hValues[177,0] = '597.00'
$hValues[178,0] = '770.00'
$hValues[179,0] = '244.00'
$hValues[180,0] = '350.00'
$hValues[181,0] = '5128.00'
$hValues[182,0] = '54246.92'
$hValues[183,0] = '1301.89'
$hValues[184,0] = '1464.00'
$hValues[185,0] = '850.00'
$hValues[186,0] = '4792.00'
$hValues[187,0] = '9005.00'
$hValues[188,0] = '270.00'
$hValues[189,0] = '22.80'
$hValues[190,0] = '1956.00'
$hValues[191,0] = '2160.00'
$hValues[192,0] = '145.25'
$hValues[193,0] = '36885.52'
$hValues[194,0] = '100.00'
$hValues[195,0] = '1085.00'
$hValues[196,0] = '825.00'
$hValues[197,0] = '1083.60'
$hValues[198,0] = '1642.55'
$hValues[199,0] = '516.00'
$hValues[200,0] = '4264.30'
$hValues[201,0] = '87.05'
$hValues[202,0] = '3475.88'
$hValues[203,0] = '3563.01'
$hValues[204,0] = '6851.02'
$hValues[205,0] = '600.00'
$hValues[206,0] = '58312.00'
$hValues[207,0] = '7200.00'
$hValues[208,0] = '651206.16'
$hValues[209,0] = '2080.00'
$hValues[210,0] = '148986.00'
$hValues[211,0] = '2323.00'
$hValues[212,0] = '1390.00'
$hValues[213,0] = '3000.00'
$hValues[214,0] = '790.00'
$hValues[215,0] = '182500.00'
$hValues[216,0] = '20000.00'
$hValues[217,0] = '83000.00'
$hValues[218,0] = '20000.00'
$hValues[219,0] = '112500.00'
$hValues[220,0] = '75000.00'
$hValues[221,0] = '223000.00'
$hValues[222,0] = '5000.00'
$hValues[223,0] = '111500.00'
$hValues[224,0] = '223000.00'
$hValues[225,0] = '257911.50'
$hValues[226,0] = '6681.19'
$hValues[227,0] = '515426.36'
$hValues[228,0] = '32530.00'
$hValues[229,0] = '10000.00'
$hValues[230,0] = '1185.62'
$hValues[231,0] = '158.34'
$hValues[232,0] = '3214.00'
$hValues[233,0] = '16000.00'
$hValues[234,0] = '6700.00'
$hValues[235,0] = '135000.00'
$hValues[236,0] = '650.00'
$hValues[237,0] = '645.00'
$hValues[238,0] = '145.00'
$hRange.Value = $hValues
$hRange.Style = "Normal"

# --- Razon social / Nombre Fantasia: normalize separators ---
$nameEdits = @(
    @('E35', 'RAMIREZ CLAUDIA. RAMIREZ CESAR Y RAMIREZ VERONICA SH'),
    @('E87', 'FERNANDEZ. MARIO HUGO'),
    @('E91', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F91', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E92', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'),
    @('E114', 'TRABICHET MARIA. VERGARA ADEL Y OTRA'),
    @('F114', 'TRABICHET MARIA. VERGARA ADEL Y OTRA'),
    @('E126', 'RICCOTTI. MARIANA EDITH'),
    @('E129', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F129', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F135', 'MERCANZINI. GASTON ARIEL'),
    @('E167', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'),
    @('E170', 'DODERA. JORGE ABELARDO'),
    @('E189', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'),
    @('E215', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH')
)
foreach ($edit in $nameEdits) {
    $cell = $ws.Range($edit[0])
    $cell.NumberFormat = "@"
    $cell.Value = $edit[1]
    $cell.Style = "Normal"
}
